{"js": "// Replace each two-digit-multiplication equation's text with its updated\n// answer. The mapping below is the exact old->new text pairs from the\n// target diff, applied in document order (each old string is unique, so a\n// literal, case-sensitive whole-match search uniquely identifies the run).\nconst replacements = [\n  [\"44\u00d763=2772\", \"86\u00d738=3268\"],\n  [\"83\u00d727=2241\", \"91\u00d745=4095\"],\n  [\"47\u00d764=3008\", \"93\u00d740=3720\"],\n  [\"57\u00d751=2907\", \"95\u00d780=7600\"],\n  [\"44\u00d739=1716\", \"79\u00d735=2765\"],\n  [\"34\u00d732=1088\", \"96\u00d777=7392\"],\n  [\"91\u00d788=8008\", \"88\u00d782=7216\"],\n  [\"41\u00d791=3731\", \"24\u00d768=1632\"],\n  [\"75\u00d771=5325\", \"48\u00d766=3168\"],\n  [\"36\u00d776=2736\", \"20\u00d724=480\"],\n  [\"32\u00d758=1856\", \"12\u00d798=1176\"],\n  [\"76\u00d782=6232\", \"28\u00d719=532\"],\n  [\"19\u00d721=399\", \"15\u00d723=345\"],\n  [\"59\u00d782=4838\", \"37\u00d744=1628\"],\n  [\"60\u00d730=1800\", \"51\u00d756=2856\"],\n  [\"46\u00d797=4462\", \"42\u00d739=1638\"],\n  [\"97\u00d751=4947\", \"51\u00d738=1938\"],\n  [\"37\u00d728=1036\", \"30\u00d728=840\"],\n  [\"36\u00d725=900\", \"90\u00d733=2970\"],\n  [\"11\u00d768=748\", \"78\u00d783=6474\"],\n  [\"41\u00d732=1312\", \"20\u00d794=1880\"],\n  [\"70\u00d737=2590\", \"13\u00d778=1014\"],\n  [\"46\u00d764=2944\", \"27\u00d763=1701\"],\n  [\"49\u00d774=3626\", \"31\u00d721=651\"],\n  [\"36\u00d795=3420\", \"37\u00d735=1295\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-multiplication equation's text with its updated\n# answer. The mapping below is the exact old->new text pairs from the\n# target diff, applied via Find/Replace (each old string is unique in the\n# document, so a whole-document ReplaceAll touches exactly one run each).\n$replacements = @(\n  @(\"44\u00d763=2772\", \"86\u00d738=3268\"),\n  @(\"83\u00d727=2241\", \"91\u00d745=4095\"),\n  @(\"47\u00d764=3008\", \"93\u00d740=3720\"),\n  @(\"57\u00d751=2907\", \"95\u00d780=7600\"),\n  @(\"44\u00d739=1716\", \"79\u00d735=2765\"),\n  @(\"34\u00d732=1088\", \"96\u00d777=7392\"),\n  @(\"91\u00d788=8008\", \"88\u00d782=7216\"),\n  @(\"41\u00d791=3731\", \"24\u00d768=1632\"),\n  @(\"75\u00d771=5325\", \"48\u00d766=3168\"),\n  @(\"36\u00d776=2736\", \"20\u00d724=480\"),\n  @(\"32\u00d758=1856\", \"12\u00d798=1176\"),\n  @(\"76\u00d782=6232\", \"28\u00d719=532\"),\n  @(\"19\u00d721=399\", \"15\u00d723=345\"),\n  @(\"59\u00d782=4838\", \"37\u00d744=1628\"),\n  @(\"60\u00d730=1800\", \"51\u00d756=2856\"),\n  @(\"46\u00d797=4462\", \"42\u00d739=1638\"),\n  @(\"97\u00d751=4947\", \"51\u00d738=1938\"),\n  @(\"37\u00d728=1036\", \"30\u00d728=840\"),\n  @(\"36\u00d725=900\", \"90\u00d733=2970\"),\n  @(\"11\u00d768=748\", \"78\u00d783=6474\"),\n  @(\"41\u00d732=1312\", \"20\u00d794=1880\"),\n  @(\"70\u00d737=2590\", \"13\u00d778=1014\"),\n  @(\"46\u00d764=2944\", \"27\u00d763=1701\"),\n  @(\"49\u00d774=3626\", \"31\u00d721=651\"),\n  @(\"36\u00d795=3420\", \"37\u00d735=1295\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
